{"js": "// The document contains a paragraph holding a Word FIELD (fldChar begin/\n// instrText.../fldChar end) whose field code is \" m:self.name \". The edit\n// converts this field into plain literal text runs reading \"{m:self.name}\"\n// (i.e. the surrounding space characters become the literal curly braces),\n// while keeping the existing run-level formatting (the orange theme color\n// applied to the \"self\" portion) intact.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Locate the paragraph that contains the field (fldChar/instrText run) by\n// loading each paragraph's field collection and checking its item count.\nconst paragraphFields = paragraphs.items.map((p) => {\n  const fields = p.fields;\n  fields.load(\"items\");\n  return fields;\n});\nawait context.sync();\n\nlet fieldParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphFields[i].items.length > 0) {\n    fieldParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (fieldParagraph) {\n  // Build the replacement paragraph OOXML: same text as the field code,\n  // split into runs the same way the original instrText runs were split,\n  // but now as literal <w:t> runs, with the leading/trailing spaces of the\n  // field code turned into the literal \"{\" and \"}\" characters, and the\n  // orange color formatting preserved on the \"self\" run.\n  const replacementOoxml =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n    '<pkg:xmlData>' +\n    '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n    '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n    '</Relationships>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n    '<w:p>' +\n    '<w:r><w:t>{</w:t></w:r>' +\n    '<w:r><w:t>m</w:t></w:r>' +\n    '<w:r><w:t>:</w:t></w:r>' +\n    '<w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>self</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">.name}</w:t></w:r>' +\n    '</w:p>' +\n    '<w:sectPr/>' +\n    '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>';\n\n  fieldParagraph.insertOoxml(replacementOoxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# The document contains a paragraph holding a Word FIELD (fldChar begin/\n# instrText.../fldChar end) whose field code is \" m:self.name \". The edit\n# converts this field into plain literal text runs reading \"{m:self.name}\"\n# (i.e. the surrounding space characters become the literal curly braces),\n# while keeping the existing run-level formatting (the orange theme color\n# applied to the \"self\" portion) intact.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that contains the field.\n$fieldParagraph = $null\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Fields.Count -gt 0) {\n        $fieldParagraph = $p\n        break\n    }\n}\n\nif ($fieldParagraph -ne $null) {\n    $range = $fieldParagraph.Range\n\n    $xml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t>{</w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>self</w:t></w:r><w:r><w:t xml:space=\"preserve\">.name}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n    $range.InsertXML($xml)\n}\n"}
